# Append some level 2 notes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 / column C: replace the numeric placeholder with the real note text.
# The text starts with '+', so (like C4 above it) Excel stores it as a
# quote-prefixed literal string rather than trying to parse it as a formula.
# Setting .Value alone does not flip the quote-prefix style bit in this
# runtime, so clone the already-quote-prefixed format from C4.
$ws.Range("C6").Value = "+根据会话状态信息，为进出数据流提供明确的允许、拒绝访问能力。"
$ws.Range("C4").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 9 / column C
$ws.Range("C9").Value = "关键网络节点监视网络攻击行为。"

# Row 10 / column C
$ws.Range("C10").Value = "关键网络节点对恶意代码进行检测和清除，并维护恶意代码防护机制的升级和更新。"

# Row 20 / column C — newly-wrapped text needs a taller row than the default.
$ws.Range("C20").Value = "划分不同的网络区域，并按照方便管理控制原则分配地址，避免将重要网络部署于边界，采用可靠技术手段隔离其他网络区域。"
$ws.Rows.Item(20).RowHeight = 41.4

# Row 21 / column C
$ws.Range("C21").Value = "网络边界和重要网络节点进行安全审计、审计覆盖每个用户，对重要的用户行为和安全事件进行审计；审计记录包括事件时间、用户、事件类型、时间是否成功与其他审计相关信息；对审计记录进行保护、定期备份，避免受到未预期的删除、修改或覆盖等。"
$ws.Rows.Item(21).RowHeight = 82.8

# Move the viewport / selection to match the author's final cursor position.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C7").Select()
